$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.628386378288269
$ws.Range("B1").Value = 3.313309192657471
$ws.Range("C1").Value = 4.037763595581055
$ws.Range("D1").Value = 1.295402765274048
$ws.Range("E1").Value = 0.7466388940811157
